$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.8
$ws.Range("H2").Value = 3.4
$ws.Range("I2").Value = 4.5
$ws.Range("K2").Value = 8
$ws.Range("N2").Value = 2.4
$ws.Range("O2").Value = 1.53
$ws.Range("S2").Value = 1.58
$ws.Range("U2").Value = 7.5
$ws.Range("W2").Value = 15
$ws.Range("Z2").Value = 7
$ws.Range("AA2").Value = 6.5
$ws.Range("AD2").Value = 9.5
$ws.Range("AE2").Value = 21
$ws.Range("S3").Value = 1.63
$ws.Range("J4").Value = 1.08
$ws.Range("K4").Value = 8
$ws.Range("L4").Value = 1.5
$ws.Range("M4").Value = 2.5
$ws.Range("S4").Value = 1.63
$ws.Range("G5").Value = 2.1
$ws.Range("I5").Value = 3.6
$ws.Range("S5").Value = 1.63
$ws.Range("U5").Value = 9
$ws.Range("AF5").Value = 13
$ws.Range("G6").Value = 3.6
$ws.Range("H6").Value = 3.4
$ws.Range("I6").Value = 2
$ws.Range("R6").Value = 1.8
$ws.Range("S6").Value = 1.8
$ws.Range("AG6").Value = 17
$ws.Range("G8").Value = 2.92
$ws.Range("H8").Value = 2.4
$ws.Range("I8").Value = 3.1
$ws.Range("J8").Value = 1.2
$ws.Range("K8").Value = 4
$ws.Range("N8").Value = 3.3
$ws.Range("Q8").Value = 1.95
$ws.Range("R8").Value = 2.37
$ws.Range("S8").Value = 1.52
$ws.Range("T8").Value = 5.5
$ws.Range("U8").Value = 12.5
$ws.Range("V8").Value = 12
$ws.Range("W8").Value = 40
$ws.Range("X8").Value = 40
$ws.Range("Y8").Value = 70
$ws.Range("Z8").Value = 4
$ws.Range("AA8").Value = 5.2
$ws.Range("AB8").Value = 21
$ws.Range("AC8").Value = 175
$ws.Range("AD8").Value = 5.9
$ws.Range("AE8").Value = 14
$ws.Range("AF8").Value = 12.5
$ws.Range("AH8").Value = 40
$ws.Range("AI8").Value = 70
$ws.Range("G9").Value = 1.66
$ws.Range("J9").Value = 1.07
$ws.Range("K9").Value = 9
$ws.Range("G11").Value = 2.6
$ws.Range("H11").Value = 3.4
$ws.Range("I11").Value = 2.45
$ws.Range("N11").Value = 1.53
$ws.Range("O11").Value = 2.18
$ws.Range("R11").Value = 1.44
$ws.Range("S11").Value = 2.4
$ws.Range("T11").Value = 13
$ws.Range("X11").Value = 18.5
$ws.Range("Y11").Value = 20
$ws.Range("Z11").Value = 14.5
$ws.Range("AA11").Value = 7
$ws.Range("AD11").Value = 11.25
$ws.Range("AE11").Value = 15
$ws.Range("AF11").Value = 9.5
$ws.Range("AH11").Value = 18
$ws.Range("AI11").Value = 21
$ws.Range("AJ11").Value = 175
$ws.Range("M13").Value = 3.7
$ws.Range("R13").Value = 1.66
$ws.Range("S13").Value = 2.09
$ws.Range("J17").Value = 1.02
$ws.Range("K17").Value = 19
$ws.Range("N17").Value = 1.5
$ws.Range("O17").Value = 2.5
$ws.Range("I18").Value = 3.8
$ws.Range("W18").Value = 17
$ws.Range("AE18").Value = 19
$ws.Range("AH18").Value = 41
$ws.Range("N20").Value = 1.93
$ws.Range("O20").Value = 1.93
$ws.Range("G21").Value = 2.8
$ws.Range("I21").Value = 2.2
$ws.Range("N21").Value = 1.7
$ws.Range("O21").Value = 2.1
$ws.Range("W21").Value = 29
$ws.Range("AE21").Value = 13
$ws.Range("J22").Value = 1.05
$ws.Range("K22").Value = 11
$ws.Range("L22").Value = 1.25
$ws.Range("M22").Value = 3.75
$ws.Range("N22").Value = 1.93
$ws.Range("O22").Value = 1.93
$ws.Range("G23").Value = 1.42
$ws.Range("J23").Value = 1.02
$ws.Range("K23").Value = 11
$ws.Range("N23").Value = 1.44
$ws.Range("O23").Value = 2.63
$ws.Range("R23").Value = 1.57
$ws.Range("S23").Value = 2.25
$ws.Range("U23").Value = 9.5
$ws.Range("W23").Value = 12
$ws.Range("X23").Value = 12
$ws.Range("Y23").Value = 19
$ws.Range("Z23").Value = 21
$ws.Range("AE23").Value = 29
$ws.Range("AH23").Value = 34
$ws.Range("AJ23").Value = 126
$ws.Range("G24").Value = 1.77
$ws.Range("G25").Value = 1.17
$ws.Range("H25").Value = 7.5
$ws.Range("I25").Value = 10
$ws.Range("J25").Value = 26
$ws.Range("K25").Value = 1.02
$ws.Range("P25").Value = 1.17
$ws.Range("R25").Value = 1.83
$ws.Range("S25").Value = 1.83
$ws.Range("V25").Value = 9.5
$ws.Range("Y25").Value = 21
$ws.Range("Z25").Value = 26
$ws.Range("AA25").Value = 17
$ws.Range("AB25").Value = 26
$ws.Range("AI25").Value = 41
$ws.Range("K26").Value = 6.3
$ws.Range("P26").Value = 1.44
$ws.Range("Q26").Value = 2.62
$ws.Range("T26").Value = 8.5
$ws.Range("U26").Value = 16
$ws.Range("V26").Value = 11.5
$ws.Range("W26").Value = 45
$ws.Range("Z26").Value = 6.3
$ws.Range("AA26").Value = 6.1
$ws.Range("AH26").Value = 19.5
$ws.Range("I27").Value = 4
$ws.Range("J27").Value = 1.07
$ws.Range("Q27").Value = 2.62
$ws.Range("T27").Value = 6.9
$ws.Range("Y27").Value = 27
$ws.Range("AC27").Value = 65
$ws.Range("AD27").Value = 11.25
